$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.9399999999999999
$ws.Range("C2").Value = 0.58
$ws.Range("F2").Value = 1.05
$ws.Range("G2").Value = 1.58

$ws.Range("B3").Value = 1.04
$ws.Range("F3").Value = 1.04

$ws.Range("B4").Value = 1.04
$ws.Range("F4").Value = 1.04
